$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete original rows 2 and 3 (Eagles / Albatros rows) so the remaining
# "Prejuveniles" rows (formerly 4,5,6) shift up to become rows 2,3,4.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Now remove the trailing row (formerly row 6, Petric) which is now row 4.
$ws.Rows.Item(4).Delete()

# Update D3 (posicion) to match target value of 2 - already 2 from the
# shifted data, but set explicitly to be safe.
$ws.Range("D3").Value = 2
